# Generate Report for Handback
# Adds a new handback row (4fceb0b3-1056-4f6b-af45-69fb45137b07.md) to the
# Overview, zh-cn and de-de sheets/tables.

$wb = $excel.ActiveWorkbook

$fileId = "4fceb0b3-1056-4f6b-af45-69fb45137b07"
$fileName = $fileId + ".md"
$pathAndName = "e2e\" + $fileName

# ---------------------------------------------------------------------------
# Sheet "Overview" (table3 / Overview)
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.ListRows.Add() | Out-Null

$wsOverview.Range("A3").Value = $fileName
$wsOverview.Range("B3").Value = $pathAndName
$wsOverview.Range("C3").Value = ".md"
$wsOverview.Range("E3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("G3").Value = "2016-09-07 07:49:19"
$wsOverview.Range("G3").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsOverview.Hyperlinks.Add($wsOverview.Range("B3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/332943c98fcb3040392d9f78ea41d0e628f9c05a/e2e/" + $fileName, "", "", $pathAndName) | Out-Null

# ---------------------------------------------------------------------------
# Sheet "zh-cn" (table1)
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$loZhCn = $wsZhCn.ListObjects.Item(1)
$loZhCn.ListRows.Add() | Out-Null

$wsZhCn.Range("A3").Value = $fileName
$wsZhCn.Range("B3").Value = ".md"
$wsZhCn.Range("C3").Value = "Handed back: in sync with en-US"
$wsZhCn.Range("D3").Value = "e2e"
$wsZhCn.Range("E3").Value = "ht"
$wsZhCn.Range("F3").Value = "True"
$wsZhCn.Range("G3").Value = $fileId + ".a4985e21b01bc897f285df48bac5c3e02c44dd69.zh-cn.xlf"
$wsZhCn.Range("H3").Value = "2016-09-07 07:49:00"
$wsZhCn.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZhCn.Range("I3").Value = $fileName
$wsZhCn.Range("J3").Value = $fileId + ".a4985e21b01bc897f285df48bac5c3e02c44dd69.zh-cn.xlf"
$wsZhCn.Range("K3").Value = "2016-09-07 07:50:00"
$wsZhCn.Range("K3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZhCn.Range("L3").Value = ""
$wsZhCn.Range("M3").Value = "True"
$wsZhCn.Range("N3").Value = ""
$wsZhCn.Range("O3").Value = "False"
$wsZhCn.Range("P3").Value = ""

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/332943c98fcb3040392d9f78ea41d0e628f9c05a/e2e/" + $fileName, "", "", $fileName) | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/003ea4d5e1295f79aeb87509b4c864d7fdadde32/e2e/" + $fileName, "", "", $fileName) | Out-Null

# ---------------------------------------------------------------------------
# Sheet "de-de" (table2)
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$loDeDe = $wsDeDe.ListObjects.Item(1)
$loDeDe.ListRows.Add() | Out-Null

$wsDeDe.Range("A3").Value = $fileName
$wsDeDe.Range("B3").Value = ".md"
$wsDeDe.Range("C3").Value = "Handed back: in sync with en-US"
$wsDeDe.Range("D3").Value = "e2e"
$wsDeDe.Range("E3").Value = "ht"
$wsDeDe.Range("F3").Value = "True"
$wsDeDe.Range("G3").Value = $fileId + ".a4985e21b01bc897f285df48bac5c3e02c44dd69.de-de.xlf"
$wsDeDe.Range("H3").Value = "2016-09-07 07:49:19"
$wsDeDe.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDeDe.Range("I3").Value = $fileName
$wsDeDe.Range("J3").Value = $fileId + ".a4985e21b01bc897f285df48bac5c3e02c44dd69.de-de.xlf"
$wsDeDe.Range("K3").Value = "2016-09-07 07:50:34"
$wsDeDe.Range("K3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDeDe.Range("L3").Value = ""
$wsDeDe.Range("M3").Value = "True"
$wsDeDe.Range("N3").Value = ""
$wsDeDe.Range("O3").Value = "False"
$wsDeDe.Range("P3").Value = ""

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/332943c98fcb3040392d9f78ea41d0e628f9c05a/e2e/" + $fileName, "", "", $fileName) | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/a76ed007fe349969a60e109407dd6fa3cfb102d3/e2e/" + $fileName, "", "", $fileName) | Out-Null

Write-Host "Handback report row added"
